$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-27 down to 10-28.
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with a new weekly record (same
# constant fields as the rest of the sheet, new date and price data).
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 44804
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = 100112043
$ws.Cells.Item(9, 7).Value = "Pepino dulce"
$ws.Cells.Item(9, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 100
$ws.Cells.Item(9, 11).Value = 15000
$ws.Cells.Item(9, 12).Value = 16000
$ws.Cells.Item(9, 13).Value = 15500
$ws.Cells.Item(9, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(9, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 16).Value = 861
$ws.Cells.Item(9, 17).Value = 18
$ws.Cells.Item(9, 18).Value = "Hortaliza"
